# Updated cryptos list (price/volume refresh) - mirrors the GitHub Actions data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.820.31"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.637.44"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.70"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3856"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.69"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.318"
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.60"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.931"
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.752"
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001300"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "1.637.24"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.39"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06925"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.32"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.850"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.45"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").Value = "23.816.50"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.427"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.855"
$ws.Range("E26").Value = "  -9.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.75"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.61"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.462"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.37"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.751"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "1.821.59"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07924"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9760"
$ws.Range("E35").Value = "  -7.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02877"
$ws.Range("E36").Value = "  -4.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.527"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2640"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.37"
$ws.Range("E39").Value = "  -8.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09057"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7448"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.415"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.18"
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.43"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6843"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.392"
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.062"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08211"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.70"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("E51").Value = "  -2.61%  "
